$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = -24
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -7
$ws.Range("F9").Value = -5
